$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "spruce"
$ws.Range("B4").Value = "birch"
$ws.Range("B5").Select()
